$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -122.5725
$ws.Range("B2").Value = -122.4701

$ws.Range("A3").Value = 37.8606
$ws.Range("B3").Value = 37.9398

$ws.Range("A4").Value = -121.9428
$ws.Range("B4").Value = -122.0446

$ws.Range("A5").Value = 38.3486
$ws.Range("B5").Value = 38.2696
